$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill G25:R25 with the same value as F25 ("BEA - Github")
$value = $ws.Range("F25").Value2
$ws.Range("G25:R25").Value = $value

# Update the view: scroll back to top and move the active selection to R25
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("R25").Select()
